$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Raluca")

# Rename the worksheet (this also updates the _FilterDatabase defined name
# and the sheet reference automatically).
$ws.Name = "Wishlist Raluca"

# Mark the item in row 25 as acquired ("Y" in column E), matching the
# pattern used by the other rows in the list.
$ws.Range("E25").Value = "Y"

# Move the active selection, matching the author's final cursor position.
$ws.Range("E32").Select()
